$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

$ws.Range("D42").Value = "python nan,inf, try~exception, locals()관련 에러 예외처리 정리"

$ws.Range("D50").Value = "위상정렬"
$ws.Range("E50").Value = "http://incredible.egloos.com/7547069"

$ws.Range("D51").Value = "블로그 이름 바꿉니다. 비스카이비전 -> 코딩재개발"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EB%B8%94%EB%A1%9C%EA%B7%B8-%EC%9D%B4%EB%A6%84-%EB%B0%94%EA%BF%89%EB%8B%88%EB%8B%A4-%EB%B9%84%EC%8A%A4%EC%B9%B4%EC%9D%B4%EB%B9%84%EC%A0%84-%EC%BD%94%EB%94%A9%EC%9E%AC%EA%B0%9C%EB%B0%9C"

$ws.Range("D52").Value = "숨은 DS"
